$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the "D/E" row (row 17) formulas: replace TotalLiabilities-based
# ratios with LongTermDebt + ShortTermDebt based ones ---
$ws.Range("E17").Value = "(LongTermDebt + ShortTermDebt) / TotalShareholdersEquity "
$ws.Range("F17").Value = "(LongTermDebtCurrent + LongTermDebtNoncurrent) / StockholdersEquity"

# --- Re-color a handful of rows to the "green" highlighted style used by
# rows 10/11 (Current Ratio / Quick Ratio) elsewhere in the sheet ---
$ws.Range("B8:F8").Font.Color = 5287936
$ws.Range("C11:F11").Font.Color = 5287936
$ws.Range("C17:F17").Font.Color = 5287936
$ws.Range("C18:F18").Font.Color = 5287936

# --- Add a new "Financial Leverage" ratio row (row 27) ---
$ws.Range("B27").Value = 23
$ws.Range("C27").Value = "Financial Leverage"
$ws.Range("D27").Value = "FinLev"
$ws.Range("E27").Value = "TotalDebt / TotalShareholdersEquity"
$ws.Range("F27").Value = "/ StockholdersEquity"

$ws.Range("B26:F26").Copy()
$ws.Range("B27:F27").PasteSpecial(-4122)

# --- Update the active selection to reflect the last-edited cell ---
$ws.Range("E17").Select() | Out-Null
